$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Refresh the cached "datetimeFigureOut" date field text on the slide
#    master and every slide layout (9/19/2017 -> 11/3/2017).
# ---------------------------------------------------------------------------
function Update-DateField($container) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $sh = $container.Shapes.Item($j)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "9/19/2017") {
                $sh.TextFrame.TextRange.Text = "11/3/2017"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DateField $layouts.Item($i)
}

# ---------------------------------------------------------------------------
# 2. Rewrite the USGS disclaimer textbox on slide 2 ("TextBox 6") with the
#    new "approved for release" wording, split across two paragraphs, and
#    resize the textbox to its new position/height.
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)

$disclaimerShape = $null
for ($j = 1; $j -le $slide2.Shapes.Count; $j++) {
    $sh = $slide2.Shapes.Item($j)
    if ($sh.Name -eq "TextBox 6") {
        $disclaimerShape = $sh
    }
}

$tr = $disclaimerShape.TextFrame.TextRange
$tr.Text = "This software has been approved for release by the U.S. Geological "
$tr.InsertAfter("Survey (USGS") | Out-Null
$tr.InsertAfter("). Although the software has been subjected to rigorous review, the ") | Out-Null
$tr.InsertAfter("USGS reserves ") | Out-Null
$tr.InsertAfter("the right to update the software as ") | Out-Null
$tr.InsertAfter("needed pursuant ") | Out-Null
$tr.InsertAfter("to further ") | Out-Null
$tr.InsertAfter("analysis and ") | Out-Null
$tr.InsertAfter("review. No warranty, expressed or implied, is made by the USGS or the ") | Out-Null
$tr.InsertAfter("U.S. Government ") | Out-Null
$tr.InsertAfter("as to the functionality of the software and related material ") | Out-Null
$tr.InsertAfter("nor shall ") | Out-Null
$tr.InsertAfter("the fact of ") | Out-Null
$tr.InsertAfter("release constitute ") | Out-Null
$tr.InsertAfter("any such warranty. Furthermore, ") | Out-Null
$tr.InsertAfter("the software ") | Out-Null
$tr.InsertAfter("is released on condition that neither the USGS nor the U.S. ") | Out-Null
$tr.InsertAfter("Government shall ") | Out-Null
$tr.InsertAfter("be held liable for any damages resulting from its authorized or") | Out-Null
$tr.InsertAfter([char]13) | Out-Null
$tr.InsertAfter("unauthorized use.") | Out-Null

# Reposition / resize the textbox (left & width stay the same; top moves up
# slightly and height grows to fit the now-longer, two-paragraph text).
$disclaimerShape.Left = 65.5220472440945
$disclaimerShape.Top = 385.02685039370084
$disclaimerShape.Width = 587.228188976378
$disclaimerShape.Height = 109.05472440944882
